# volume control fix and battery level checker
#
# - Row 17 ("Improve audio quality") status flips from "Open" to "Ongoing"
#   (this is the "volume control fix" referenced in the commit message --
#   Volume Control/audio quality work moves from Open to Ongoing).
# - Three new todo rows are appended for the battery level checker and
#   related items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Improve audio quality: Open -> Ongoing
$ws.Range("C17").Value = "Ongoing"

# New rows appended after the existing last row (43)
$ws.Range("B44").Value = "Battery level check with led bar"
$ws.Range("C44").Value = "Done"

$ws.Range("B45").Value = "Improve streaming performance"
$ws.Range("C45").Value = "Open"

$ws.Range("B46").Value = "App support for intercom moving to new WiFi network"
$ws.Range("C46").Value = "Open"

# Match the author's final selection/scroll position as closely as the
# object model allows.
$ws.Range("C48").Select() | Out-Null
